$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Column A labels (rows 1-31) and Column B values (rows 1-31).
# Rows with $null value become blank numeric cells (no <v> content),
# matching the "Lương cơ bản" rows that have not yet been computed.

$labels = @(
    "Danh mục",
    "Ngày công",
    "Phụ cấp",
    "Lương cơ bản tại CẦN THƠ",
    "Chiết khấu sale chính tại CẦN THƠ",
    "Chiết khấu sale phụ tại CẦN THƠ",
    "Đơn 1 bác sĩ tại CẦN THƠ",
    "Đơn 2 bác sĩ tại CẦN THƠ",
    "Công phụ phẫu 1 tại CẦN THƠ",
    "Công phụ phẫu 2 tại CẦN THƠ",
    "Ứng lương tại CẦN THƠ",
    "Lương cơ bản tại LONG XUYÊN",
    "Chiết khấu sale chính tại LONG XUYÊN",
    "Chiết khấu sale phụ tại LONG XUYÊN",
    "Đơn 1 bác sĩ tại LONG XUYÊN",
    "Đơn 2 bác sĩ tại LONG XUYÊN",
    "Công phụ phẫu 1 tại LONG XUYÊN",
    "Công phụ phẫu 2 tại LONG XUYÊN",
    "Ứng lương tại LONG XUYÊN",
    "Lương cơ bản tại SÓC TRĂNG",
    "Chiết khấu sale chính tại SÓC TRĂNG",
    "Chiết khấu sale phụ tại SÓC TRĂNG",
    "Đơn 1 bác sĩ tại SÓC TRĂNG",
    "Đơn 2 bác sĩ tại SÓC TRĂNG",
    "Công phụ phẫu 1 tại SÓC TRĂNG",
    "Công phụ phẫu 2 tại SÓC TRĂNG",
    "Ứng lương tại SÓC TRĂNG",
    "Tổng lương tại CẦN THƠ",
    "Tổng lương tại LONG XUYÊN",
    "Tổng lương tại SÓC TRĂNG",
    "Tổng lương"
)

$values = @(
    3,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    $null,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    $null,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    if ($null -eq $values[$i]) {
        $ws.Cells.Item($row, 2).Value = $null
    } else {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}
